$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the selected cell (cosmetic change reflecting last selection before save)
$ws.Range("F29").Select()

# Update the raw input values that drove the recalculated systematic error table
$ws.Range("B25:E25").Value = 10

$ws.Range("B28:D28").Value = 45
$ws.Range("E28").Value = 58

$wb.Save()
